# Auto-generated edit script applying profit/price recalculations
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1459825.8
$ws.Range("I70").Value = 5103540
$ws.Range("J70").Value = 2340
$ws.Range("K70").Value = 15310620
$ws.Range("L70").Value = 7020
$ws.Range("M70").Value = -15310350
$ws.Range("N70").Value = -7560
$ws.Range("H73").Value = 1459825.8
$ws.Range("I73").Value = 5103540
$ws.Range("J73").Value = 2340
$ws.Range("K73").Value = 15310620
$ws.Range("L73").Value = 7020
$ws.Range("M73").Value = -15309684
$ws.Range("N73").Value = -8892
$ws.Range("H76").Value = 10547030
$ws.Range("I76").Value = 2900
$ws.Range("J76").Value = 12053334
$ws.Range("K76").Value = 2900
$ws.Range("L76").Value = 12053334
$ws.Range("M76").Value = -2585
$ws.Range("N76").Value = -12053964
$ws.Range("H79").Value = 10547030
$ws.Range("I79").Value = 2900
$ws.Range("J79").Value = 12053334
$ws.Range("K79").Value = 2900
$ws.Range("L79").Value = 12053334
$ws.Range("M79").Value = -1808
$ws.Range("N79").Value = -12055518
$ws.Range("H98").Value = 3603.3333
$ws.Range("I98").Value = 3967.8125
$ws.Range("J98").Value = 687.5
$ws.Range("K98").Value = 3967.8125
$ws.Range("L98").Value = 687.5
$ws.Range("M98").Value = -2469.8125
$ws.Range("N98").Value = -3683.5
$ws.Range("H118").Value = 270
$ws.Range("I118").Value = 270
$ws.Range("K118").Value = 810
$ws.Range("M118").Value = 847
$ws.Range("H122").Value = 3603.3333
$ws.Range("I122").Value = 3967.8125
$ws.Range("J122").Value = 687.5
$ws.Range("K122").Value = 11903.4375
$ws.Range("L122").Value = 2062.5
$ws.Range("M122").Value = -9453.4375
$ws.Range("N122").Value = -6962.5
$ws.Range("H138").Value = 2423.3157
$ws.Range("J138").Value = 2339.4614
$ws.Range("L138").Value = 7018.3842
$ws.Range("N138").Value = -17298.3842

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 4799.6665
$ws.Range("I10").Value = 4799.6665
$ws.Range("K10").Value = 4799.6665
$ws.Range("M10").Value = -4629.6665
$ws.Range("H12").Value = 620
$ws.Range("I12").Value = 525
$ws.Range("K12").Value = 525
$ws.Range("M12").Value = -352
$ws.Range("H14").Value = 699
$ws.Range("I14").Value = 394.5
$ws.Range("J14").Value = 1003.5
$ws.Range("K14").Value = 394.5
$ws.Range("L14").Value = 1003.5
$ws.Range("M14").Value = -219.5
$ws.Range("N14").Value = -1353.5
$ws.Range("H32").Value = 4728.9766
$ws.Range("I32").Value = 3755.077
$ws.Range("K32").Value = 3755.077
$ws.Range("M32").Value = -3468.077
$ws.Range("H45").Value = 2499.3333
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 2499.3333
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 2499.3333
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -3253.3333
$ws.Range("H74").Value = 30305166
$ws.Range("I74").Value = 43480120
$ws.Range("K74").Value = 43480120
$ws.Range("M74").Value = -43479246
$ws.Range("H77").Value = 30305166
$ws.Range("I77").Value = 43480120
$ws.Range("K77").Value = 217400600
$ws.Range("M77").Value = -217396232
$ws.Range("H132").Value = 6546.7144
$ws.Range("I132").Value = 6478.52
$ws.Range("J132").Value = 7115
$ws.Range("K132").Value = 19435.56
$ws.Range("L132").Value = 21345
$ws.Range("M132").Value = -16905.56
$ws.Range("N132").Value = -26405

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20039.857
$ws.Range("I20").Value = 37859.145
$ws.Range("J20").Value = 2220.5715
$ws.Range("K20").Value = 37859.145
$ws.Range("L20").Value = 2220.5715
$ws.Range("M20").Value = -37612.145
$ws.Range("N20").Value = -2714.5715
$ws.Range("H80").Value = 270.1
$ws.Range("I80").Value = 502
$ws.Range("J80").Value = 212.125
$ws.Range("K80").Value = 502
$ws.Range("L80").Value = 212.125
$ws.Range("M80").Value = 496
$ws.Range("N80").Value = -2208.125
$ws.Range("H83").Value = 270.1
$ws.Range("I83").Value = 502
$ws.Range("J83").Value = 212.125
$ws.Range("K83").Value = 2510
$ws.Range("L83").Value = 1060.625
$ws.Range("M83").Value = 2482
$ws.Range("N83").Value = -11044.625
$ws.Range("H107").Value = 1736
$ws.Range("I107").Value = 1143
$ws.Range("K107").Value = 1143
$ws.Range("M107").Value = 777

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1899.4166
$ws.Range("I15").Value = 2261.625
$ws.Range("K15").Value = 2261.625
$ws.Range("M15").Value = -2091.625
$ws.Range("H16").Value = 860
$ws.Range("I16").Value = 860
$ws.Range("K16").Value = 860
$ws.Range("M16").Value = -573
$ws.Range("H29").Value = 18332.334
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 18332.334
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 18332.334
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -18918.334
$ws.Range("H113").Value = 860
$ws.Range("I113").Value = 860
$ws.Range("K113").Value = 860
$ws.Range("M113").Value = 1310

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 791.4167
$ws.Range("I14").Value = 791.4167
$ws.Range("K14").Value = 2374.2501
$ws.Range("M14").Value = -2201.2501
$ws.Range("H24").Value = 30
$ws.Range("I24").Value = 25
$ws.Range("J24").Value = 33.333332
$ws.Range("K24").Value = 75
$ws.Range("L24").Value = 99.999996
$ws.Range("M24").Value = 155
$ws.Range("N24").Value = -559.999996
$ws.Range("H33").Value = 115.64286
$ws.Range("I33").Value = 120.3
$ws.Range("J33").Value = 104
$ws.Range("K33").Value = 721.8
$ws.Range("L33").Value = 624
$ws.Range("M33").Value = -438.8
$ws.Range("N33").Value = -1190
$ws.Range("H38").Value = 188.35
$ws.Range("I38").Value = 156.36363
$ws.Range("J38").Value = 227.44444
$ws.Range("K38").Value = 469.09089
$ws.Range("L38").Value = 682.33332
$ws.Range("M38").Value = -122.09089
$ws.Range("N38").Value = -1376.33332
$ws.Range("H69").Value = 1886
$ws.Range("I69").Value = 1886
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 5658
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -4847
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 1886
$ws.Range("I72").Value = 1886
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 16974
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -12918
$ws.Range("N72").ClearContents()
$ws.Range("H80").Value = 6235.476
$ws.Range("I80").Value = 5998
$ws.Range("J80").Value = 6829.1665
$ws.Range("K80").Value = 17994
$ws.Range("L80").Value = 20487.4995
$ws.Range("M80").Value = -17058
$ws.Range("N80").Value = -22359.4995
$ws.Range("H83").Value = 6235.476
$ws.Range("I83").Value = 5998
$ws.Range("J83").Value = 6829.1665
$ws.Range("K83").Value = 53982
$ws.Range("L83").Value = 61462.4985
$ws.Range("M83").Value = -49302
$ws.Range("N83").Value = -70822.4985
$ws.Range("H141").Value = 5659.8423
$ws.Range("I141").Value = 5619.8823
$ws.Range("K141").Value = 16859.6469
$ws.Range("M141").Value = -11679.6469

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9715.526
$ws.Range("I70").Value = 9281.182000000001
$ws.Range("K70").Value = 9281.182000000001
$ws.Range("M70").Value = -9011.182000000001
$ws.Range("H73").Value = 9715.526
$ws.Range("I73").Value = 9281.182000000001
$ws.Range("K73").Value = 9281.182000000001
$ws.Range("M73").Value = -8345.182000000001
$ws.Range("H102").Value = 1224.6865
$ws.Range("I102").Value = 723.9464
$ws.Range("K102").Value = 723.9464
$ws.Range("M102").Value = 898.0536
$ws.Range("H111").Value = 59293
$ws.Range("J111").Value = 59293
$ws.Range("L111").Value = 59293
$ws.Range("N111").Value = -65427
$ws.Range("H132").Value = 2678.8914
$ws.Range("I132").Value = 2810.3784
$ws.Range("J132").Value = 2138.3333
$ws.Range("K132").Value = 8431.135200000001
$ws.Range("L132").Value = 6414.999899999999
$ws.Range("M132").Value = -5901.135200000001
$ws.Range("N132").Value = -11474.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2920
$ws.Range("I7").Value = 2446.4092
$ws.Range("K7").Value = 2446.4092
$ws.Range("M7").Value = -2334.4092
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 802
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 802
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1142
$ws.Range("H26").Value = 9833
$ws.Range("I26").Value = 9833
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 9833
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -9538
$ws.Range("N26").ClearContents()
$ws.Range("H40").Value = 1767408.2
$ws.Range("I40").Value = 1877533.8
$ws.Range("J40").Value = 5400
$ws.Range("K40").Value = 1877533.8
$ws.Range("L40").Value = 5400
$ws.Range("M40").Value = -1877397.8
$ws.Range("N40").Value = -5672
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H126").Value = 2920
$ws.Range("I126").Value = 2446.4092
$ws.Range("K126").Value = 7339.2276
$ws.Range("M126").Value = -4869.2276

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1707923.1
$ws.Range("I62").Value = 2982428
$ws.Range("K62").Value = 2982428
$ws.Range("M62").Value = -2981804
$ws.Range("H65").Value = 1707923.1
$ws.Range("I65").Value = 2982428
$ws.Range("K65").Value = 14912140
$ws.Range("M65").Value = -14909020
$ws.Range("H107").Value = 1007.5714
$ws.Range("I107").Value = 915
$ws.Range("K107").Value = 2745
$ws.Range("M107").Value = -825
$ws.Range("H113").Value = 2291
$ws.Range("I113").Value = 777.8570999999999
$ws.Range("J113").Value = 3350.2
$ws.Range("K113").Value = 2333.5713
$ws.Range("L113").Value = 10050.6
$ws.Range("M113").Value = -163.5712999999996
$ws.Range("N113").Value = -14390.6
$ws.Range("H132").Value = 19892
$ws.Range("J132").Value = 7224
$ws.Range("L132").Value = 21672
$ws.Range("N132").Value = -26732

